$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column N (weights) updates for rows 2-21 (some rows unchanged and intentionally omitted: N9, N12)
$ws.Range("N2").Value = 0.001000000000000012
$ws.Range("N3").Value = 0.1341018387334771
$ws.Range("N4").Value = 0.001
$ws.Range("N5").Value = 0.15
$ws.Range("N6").Value = 0.001000000000000024
$ws.Range("N7").Value = 0.1309300933923813
$ws.Range("N8").Value = 0.004282729605298713
$ws.Range("N10").Value = 0.007392408374097604
$ws.Range("N11").Value = 0.15
$ws.Range("N13").Value = 0.0121809121607749
$ws.Range("N14").Value = 0.07481878778029248
$ws.Range("N15").Value = 0.15
$ws.Range("N16").Value = 0.001000000000000002
$ws.Range("N17").Value = 0.00100000000000001
$ws.Range("N18").Value = 0.07431405926470459
$ws.Range("N19").Value = 0.001000000000000013
$ws.Range("N20").Value = 0.001000000000000005
$ws.Range("N21").Value = 0.1029791706889733

# Row 22 ("Portfolio return ln") updates for columns B..M
$ws.Range("B22").Value = -0.0116739433266411
$ws.Range("C22").Value = 0.0268544916715576
$ws.Range("D22").Value = 0.01885718343157542
$ws.Range("E22").Value = 0.03158375319533705
$ws.Range("F22").Value = -0.008889012310312441
$ws.Range("G22").Value = -0.01383340334412321
$ws.Range("H22").Value = 0.02138974457299077
$ws.Range("I22").Value = 0.03290845601103991
$ws.Range("J22").Value = 0.03431537956617287
$ws.Range("K22").Value = -0.07049160495096506
$ws.Range("L22").Value = -0.03443219069318885
$ws.Range("M22").Value = 0.02814705398979437

# Row 23 ("Portfolio return") updates for columns B..N
$ws.Range("B23").Value = 0.9883939327657866
$ws.Range("C23").Value = 1.027218323067348
$ws.Range("D23").Value = 1.019036102985162
$ws.Range("E23").Value = 1.032087812628083
$ws.Range("F23").Value = 0.9911503781590867
$ws.Range("G23").Value = 0.9862618385011598
$ws.Range("H23").Value = 1.02162014496229
$ws.Range("I23").Value = 1.033455928232752
$ws.Range("J23").Value = 1.034910945030068
$ws.Range("K23").Value = 0.9319355631377104
$ws.Range("L23").Value = 0.9661538516882606
$ws.Range("M23").Value = 1.028546925229916
$ws.Range("N23").Value = 1.132674325814483
